$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.143.56'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.844.62'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.48'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6972'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9986'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07705'
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3058'
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.57'
$ws.Range('E10').Value = '  -1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07821'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '92.88'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.839.22'
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.121'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6844'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.638'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008292'
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.120.79'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.83'
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.079.35'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.75'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9988'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.481'
$ws.Range('E23').Value = '  -1.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9990'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1508'
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.19'
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.815'
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.23'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.542'
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.231'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.178'
$ws.Range('E31').Value = '  -1.71%  '
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05113'
$ws.Range('E33').Value = '  -2.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7851'
$ws.Range('E34').Value = '  +3.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.863'
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.146'
$ws.Range('E36').Value = '  -2.59%  '
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.313.73'
$ws.Range('E38').Value = '  +7.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01867'
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.706'
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('E41').Value = '  +5.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.074'
$ws.Range('E42').Value = '  +5.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '107.70'
$ws.Range('E43').Value = '  -2.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9987'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.723'
$ws.Range('E45').Value = '  +2.29%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000123'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5176'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.981.67'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '64.16'
$ws.Range('E49').Value = '  -2.20%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.761'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.983'
$ws.Range('E51').Value = '  -0.77%  '
